# Musician Finder Specification Document - "Activity Area Building - Still
# Building + CSV file for it" edit.
#
# 1) The cached text of the auto-updating "today" date field
#    (datetimeFigureOut) rolls from 9/25/2020 -> 9/26/2020 on the slide
#    master and on every slide layout.
# 2) On slide 3 ("Activity Area"), the "Connect To Google Map" rectangle
#    grows taller and is relabeled "Implemented Using JSON via Dropdown",
#    and the elbow connector glued to it is repositioned/resized to match.

$p = $ppt.ActivePresentation

function Set-DateFieldText {
    param($shapes, $newText)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# -- Slide master date field --
Set-DateFieldText $p.SlideMaster.Shapes "9/26/2020"

# -- Every slide layout's date field --
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Set-DateFieldText $layouts.Item($i).Shapes "9/26/2020"
}

# -- Slide 3: "Activity Area" shapes --
$s3 = $p.Slides.Item(3)

$rect = $s3.Shapes.Item("Rectangle 7")
$rect.Height = 65.1904
$rect.TextFrame.TextRange.Text = "Implemented Using JSON via Dropdown"

$conn = $s3.Shapes.Item("Connector: Elbow 9")
$conn.Top = 235.2318
$conn.Height = 43.50071
